# Fix boolean + float cell values across the three worksheets, matching the
# FsSpreadsheet export fix: booleans are written as proper Excel boolean
# values (serialized as 1/0) instead of literal "True"/"False" text, and the
# numeric "Numbers" sample value in the last data row becomes a float.

$wb = $excel.ActiveWorkbook

# Sheet "WithTable"
$ws1 = $wb.Worksheets.Item("WithTable")
$ws1.Range("D2").Value = $true
$ws1.Range("D3").Value = $false
$ws1.Range("D4").Value = $true
$ws1.Range("D5").Value = $false
$ws1.Range("A5").Value = 4.269

# Sheet "Tableless"
$ws2 = $wb.Worksheets.Item("Tableless")
$ws2.Range("D2").Value = $true
$ws2.Range("D3").Value = $false
$ws2.Range("D4").Value = $true
$ws2.Range("D5").Value = $false
$ws2.Range("A5").Value = 4.269

# Sheet "WithTable_Duplicate"
$ws3 = $wb.Worksheets.Item("WithTable_Duplicate")
$ws3.Range("E5").Value = $true
$ws3.Range("E6").Value = $false
$ws3.Range("E7").Value = $true
$ws3.Range("E8").Value = $false
$ws3.Range("B8").Value = 4.269
